# Add a new "胎號" (litter number) column to the 分娩資料 (farrowing data)
# sheet: header in M1, and a RANDBETWEEN(1000,2000) formula filled down
# M2:M63 (entered as two operations so Excel records M3:M63 as one shared
# formula group and M2 as a separate, non-shared formula cell — matching
# how the column was actually built up).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # 分娩資料

# Header
$ws.Range("M1").Value = "胎號"

# Fill the formula: first the bulk M3:M63 range (creates the shared-formula
# block), then M2 on its own (separate formula cell).
$ws.Range("M3:M63").Formula = "=RANDBETWEEN(1000, 2000)"
$ws.Range("M2").Formula = "=RANDBETWEEN(1000, 2000)"

# Match the selection state left behind after filling the column.
$ws.Range("M2:M63").Select()

# The other two sheets (配種資料, 基本資料) had their per-row cached height
# overrides dropped when the workbook was normalised by the newer Excel
# build; AutoFit reproduces that (it clears the redundant ht= when it
# matches the sheet's default row height).
$ws2 = $wb.Worksheets.Item(2)  # 配種資料
$ws2.Range("A1:P58").Rows.AutoFit()

$ws3 = $wb.Worksheets.Item(3)  # 基本資料
$ws3.Range("A1:Q101").Rows.AutoFit()
